# MCTK2 interface prototype - date bump + "verify" relabeling edit
#
# Summary of the change (see commit message / xml diff):
#   1. The "datetimeFigureOut" footer field on the slide master and every
#      slide layout is bumped from 2017/11/28 to 2017/11/30.
#   2. On slide 1 and slide 2, the "性质编辑器" ("property editor") label is
#      renamed to "验证" ("verify"), and the neighbouring "反例展示"
#      ("counter-example display") label shape is removed entirely.
#   3. On slide 3, the layout is reworked: the "菜单栏" title becomes
#      "反例界面", the "快捷图标栏" strip is removed, the two big panels
#      ("反例图" / "结点状态..." ) grow upward to fill the reclaimed space,
#      and the leftover "模型编辑器" / "性质编辑器" / "反例展示" labels are
#      removed.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the date footer field everywhere it appears: the slide master
#    plus all eleven slide layouts.
# ---------------------------------------------------------------------
$oldDate = "2017/11/28"
$newDate = "2017/11/30"

for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    $sh = $p.SlideMaster.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1: "性质编辑器" -> "验证", drop the "反例展示" shape.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item("矩形 8").TextFrame.TextRange.Text = "验证"
$s1.Shapes.Item("矩形 9").Delete()

# ---------------------------------------------------------------------
# 3) Slide 2: same relabeling as slide 1.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item("矩形 22").TextFrame.TextRange.Text = "验证"
$s2.Shapes.Item("矩形 23").Delete()

# ---------------------------------------------------------------------
# 4) Slide 3: retitle, drop the toolbar strip and the stray labels, and
#    grow the two remaining big panels upward to reclaim that space.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# "菜单栏" -> "反例界面"
$s3.Shapes.Item("矩形 3").TextFrame.TextRange.Text = "反例界面"

# Drop the "快捷图标栏" strip.
$s3.Shapes.Item("矩形 4").Delete()

# "反例图" panel grows up to where the strip used to start.
$panelLeft = $s3.Shapes.Item("矩形 5")
$panelLeft.Top = 73.24133858
$panelLeft.Height = 439.03450819

# "结点状态 / 信息的文本输出" panel grows up the same way.
$panelRight = $s3.Shapes.Item("矩形 10")
$panelRight.Top = 73.24133858
$panelRight.Height = 439.03440945

# Remove the leftover row of labels underneath the old strip.
$s3.Shapes.Item("矩形 11").Delete()
$s3.Shapes.Item("矩形 12").Delete()
$s3.Shapes.Item("矩形 13").Delete()
